$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.190.33"
$ws.Range("E2").Value = "  -0.17%  "
$ws.Range("D3").Value = "2.635.53"
$ws.Range("E3").Value = "  -0.46%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "527.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.66%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.59"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.76%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.569"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.67"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.64%  "
$ws.Range("E10").Value = "  +1.53%  "
$ws.Range("E11").Value = "  +0.82%  "
$ws.Range("E12").Value = "  +1.08%  "
$ws.Range("D13").Value = "3.104.98"
$ws.Range("E13").Value = "  -0.43%  "
$ws.Range("D14").Value = "59.151.40"
$ws.Range("E14").Value = "  -0.28%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.02"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.12%  "
$ws.Range("E16").Value = "  +0.82%  "
$ws.Range("D17").Value = "2.647.01"
$ws.Range("E17").Value = "  -0.51%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "341.50"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.72%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.45"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.56%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.58"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.65%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.33"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.85%  "
$ws.Range("E22").Value = "  -0.25%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.34"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.60%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.418"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.75%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.167"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.06%  "
$ws.Range("E26").Value = "  -0.63%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.25"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.60%  "
$ws.Range("D28").Value = "0.0₃0799"
$ws.Range("E28").Value = "  -0.33%  "
$ws.Range("B29").Value = "USDe"
$ws.Range("C29").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.998"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("B30").Value = "Aptos"
$ws.Range("C30").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.43"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.64%  "
$ws.Range("E31").Value = "  +2.50%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.92"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.07%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "150.12"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.38%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.20"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.31%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.20"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.51%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.908"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.87%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.869"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.12%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.49"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.31%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.63"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.74%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.65"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.21%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.997"
$ws.Range("D41").Style = "Normal"
$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0976"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.20%  "
$ws.Range("B43").Value = "Bittensor"
$ws.Range("C43").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "271.97"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.32%  "
$ws.Range("B44").Value = "Mantle"
$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.601"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.38%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "19.39"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.87%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0538"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.56%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.66"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.43%  "
$ws.Range("D48").Value = "2.045.24"
$ws.Range("E48").Value = "  -0.72%  "
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0230"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.44%  "
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.69"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.89%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "18.85"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.36%  "
